# Update "想去人数" (F column) values across sheets to the newly generated
# output numbers, as published to gh-pages at commit 456a3b4.

$wb = $excel.ActiveWorkbook

# --- Sheet "展览" ---
$ws1 = $wb.Worksheets.Item("展览")
$sheet1Updates = @{
    4  = 10069
    5  = 722
    6  = 11
    7  = 188
    8  = 392
    9  = 409
    10 = 455
    11 = 234
    12 = 511
    13 = 12662
    14 = 12662
    15 = 42
    16 = 318
    20 = 47
    23 = 268
    28 = 2752
    31 = 2111
    33 = 126
    36 = 1080
    37 = 4253
    39 = 810
    41 = 63
    42 = 1362
    44 = 44
    45 = 484
    46 = 645
    47 = 77
    48 = 283
    49 = 118
    50 = 167
}
foreach ($row in $sheet1Updates.Keys) {
    $ws1.Range("F$row").Value = $sheet1Updates[$row]
}

# --- Sheet "演出" ---
$ws2 = $wb.Worksheets.Item("演出")
$sheet2Updates = @{
    4 = 8
    5 = 49
    8 = 60
}
foreach ($row in $sheet2Updates.Keys) {
    $ws2.Range("F$row").Value = $sheet2Updates[$row]
}

# --- Sheet "全部类型" ---
$ws4 = $wb.Worksheets.Item("全部类型")
$sheet4Updates = @{
    5  = 10069
    6  = 722
    7  = 11
    8  = 188
    9  = 392
    10 = 234
    11 = 12662
    12 = 318
    17 = 268
    24 = 2111
    26 = 126
    30 = 1080
    34 = 4253
    36 = 810
    38 = 63
    40 = 1362
    43 = 44
    44 = 484
    46 = 645
    47 = 77
    48 = 283
    49 = 118
    50 = 167
}
foreach ($row in $sheet4Updates.Keys) {
    $ws4.Range("F$row").Value = $sheet4Updates[$row]
}
